$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "20869_1"
$ws.Range("B2").Value = "https://oleks-netizen.github.io/product-images/20869_1/1.jpg,https://oleks-netizen.github.io/product-images/20869_1/2.jpg,https://oleks-netizen.github.io/product-images/20869_1/2.jpg,https://oleks-netizen.github.io/product-images/20869_1/10.jpg,https://oleks-netizen.github.io/product-images/20869_1/11.jpg,https://oleks-netizen.github.io/product-images/20869_1/3.jpg,https://oleks-netizen.github.io/product-images/20869_1/4.jpg,https://oleks-netizen.github.io/product-images/20869_1/5.jpg,https://oleks-netizen.github.io/product-images/20869_1/6.jpg,https://oleks-netizen.github.io/product-images/20869_1/7.jpg,https://oleks-netizen.github.io/product-images/20869_1/8.jpg"
$ws.Range("C2").Value = 11

$ws.Range("A3").Value = "20871_1"
$ws.Range("B3").Value = "https://oleks-netizen.github.io/product-images/20871_1/1.jpg,https://oleks-netizen.github.io/product-images/20871_1/2.jpg,https://oleks-netizen.github.io/product-images/20871_1/10.jpg,https://oleks-netizen.github.io/product-images/20871_1/12.jpg,https://oleks-netizen.github.io/product-images/20871_1/11.jpg,https://oleks-netizen.github.io/product-images/20871_1/3.jpg,https://oleks-netizen.github.io/product-images/20871_1/4.jpg,https://oleks-netizen.github.io/product-images/20871_1/5.jpg,https://oleks-netizen.github.io/product-images/20871_1/6.jpg,https://oleks-netizen.github.io/product-images/20871_1/7.jpg,https://oleks-netizen.github.io/product-images/20871_1/8.jpg"
$ws.Range("C3").Value = 11

$ws.Range("A4").Value = "'23148"
$ws.Range("B4").Value = "https://oleks-netizen.github.io/product-images/23148/1.jpg,https://oleks-netizen.github.io/product-images/23148/2.jpg,https://oleks-netizen.github.io/product-images/23148/9.jpg,https://oleks-netizen.github.io/product-images/23148/12.jpg,https://oleks-netizen.github.io/product-images/23148/11.jpg,https://oleks-netizen.github.io/product-images/23148/3.jpg,https://oleks-netizen.github.io/product-images/23148/4.jpg,https://oleks-netizen.github.io/product-images/23148/5.jpg,https://oleks-netizen.github.io/product-images/23148/6.jpg,https://oleks-netizen.github.io/product-images/23148/7.jpg,https://oleks-netizen.github.io/product-images/23148/8.jpg"
$ws.Range("C4").Value = 11

$ws.Range("A5").Value = "'23150"
$ws.Range("B5").Value = "https://oleks-netizen.github.io/product-images/23150/1.jpg,https://oleks-netizen.github.io/product-images/23150/2.jpg,https://oleks-netizen.github.io/product-images/23150/11.jpg,https://oleks-netizen.github.io/product-images/23150/3.jpg,https://oleks-netizen.github.io/product-images/23150/4.jpg,https://oleks-netizen.github.io/product-images/23150/5.jpg,https://oleks-netizen.github.io/product-images/23150/6.jpg,https://oleks-netizen.github.io/product-images/23150/7.jpg,https://oleks-netizen.github.io/product-images/23150/8.jpg"
$ws.Range("C5").Value = 9

$ws.Range("A6").Value = "'23152"
$ws.Range("B6").Value = "https://oleks-netizen.github.io/product-images/23152/1.jpg,https://oleks-netizen.github.io/product-images/23152/2.jpg,https://oleks-netizen.github.io/product-images/23152/9.jpg,https://oleks-netizen.github.io/product-images/23152/11.jpg,https://oleks-netizen.github.io/product-images/23152/3.jpg,https://oleks-netizen.github.io/product-images/23152/4.jpg,https://oleks-netizen.github.io/product-images/23152/5.jpg,https://oleks-netizen.github.io/product-images/23152/6.jpg,https://oleks-netizen.github.io/product-images/23152/7.jpg,https://oleks-netizen.github.io/product-images/23152/8.jpg"
$ws.Range("C6").Value = 10

$ws.Range("A7").Value = "51411030m"
$ws.Range("B7").Value = "https://oleks-netizen.github.io/product-images/51411030m/1.jpg,https://oleks-netizen.github.io/product-images/51411030m/3.jpg,https://oleks-netizen.github.io/product-images/51411030m/3.jpg,https://oleks-netizen.github.io/product-images/51411030m/7.jpg,https://oleks-netizen.github.io/product-images/51411030m/7.jpg,https://oleks-netizen.github.io/product-images/51411030m/4.jpg,https://oleks-netizen.github.io/product-images/51411030m/5.jpg,https://oleks-netizen.github.io/product-images/51411030m/6.jpg"
$ws.Range("C7").Value = 8

$ws.Range("A8").Value = "C1HS1890bl-black"
$ws.Range("B8").Value = "https://oleks-netizen.github.io/product-images/C1HS1890bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1HS1890bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1HS1890bl-black/6.jpg,https://oleks-netizen.github.io/product-images/C1HS1890bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1HS1890bl-black/5.jpg,https://oleks-netizen.github.io/product-images/C1HS1890bl-black/5.jpg"
$ws.Range("C8").Value = 6

$ws.Range("A9").Value = "C1HSSA0546bl-black"
$ws.Range("B9").Value = "https://oleks-netizen.github.io/product-images/C1HSSA0546bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546bl-black/5.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546bl-black/4.jpg"
$ws.Range("C9").Value = 5

$ws.Range("A10").Value = "C1HSSA0546gr-green"
$ws.Range("B10").Value = "https://oleks-netizen.github.io/product-images/C1HSSA0546gr-green/1.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546gr-green/5.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546gr-green/5.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546gr-green/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA0546gr-green/4.jpg"
$ws.Range("C10").Value = 5

$ws.Range("A11").Value = "C1HSSA6020bl-black"
$ws.Range("B11").Value = "https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/6.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/6.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/5.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020bl-black/5.jpg"
$ws.Range("C11").Value = 8

$ws.Range("A12").Value = "C1HSSA6020gr-green"
$ws.Range("B12").Value = "https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/1.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/2.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/2.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/3.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/4.jpg,https://oleks-netizen.github.io/product-images/C1HSSA6020gr-green/6.jpg"
$ws.Range("C12").Value = 7

$ws.Range("A13").Value = "C1SA6019bl-black"
$ws.Range("B13").Value = "https://oleks-netizen.github.io/product-images/C1SA6019bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1SA6019bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1SA6019bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1SA6019bl-black/4.jpg,https://oleks-netizen.github.io/product-images/C1SA6019bl-black/6.jpg,https://oleks-netizen.github.io/product-images/C1SA6019bl-black/3.jpg"
$ws.Range("C13").Value = 6

$ws.Range("A14").Value = "C1SA9208bl-black"
$ws.Range("B14").Value = "https://oleks-netizen.github.io/product-images/C1SA9208bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1SA9208bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1SA9208bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1SA9208bl-black/5.jpg"
$ws.Range("C14").Value = 4

$ws.Range("A15").Value = "C1YM1725bl-black"
$ws.Range("B15").Value = "https://oleks-netizen.github.io/product-images/C1YM1725bl-black/1.jpg,https://oleks-netizen.github.io/product-images/C1YM1725bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1YM1725bl-black/2.jpg,https://oleks-netizen.github.io/product-images/C1YM1725bl-black/3.jpg,https://oleks-netizen.github.io/product-images/C1YM1725bl-black/6.jpg,https://oleks-netizen.github.io/product-images/C1YM1725bl-black/4.jpg"
$ws.Range("C15").Value = 6

$ws.Range("A16").Value = "FA-7122-4x"
$ws.Range("B16").Value = "https://oleks-netizen.github.io/product-images/FA-7122-4x/1.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/2.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/3.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/4.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/10.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/6.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/7.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/8.jpg,https://oleks-netizen.github.io/product-images/FA-7122-4x/9.jpg"
$ws.Range("C16").Value = 9

$ws.Range("A17").Value = "JD7348A"
$ws.Range("B17").Value = "https://oleks-netizen.github.io/product-images/JD7348A/2.jpg,https://oleks-netizen.github.io/product-images/JD7348A/3.jpg,https://oleks-netizen.github.io/product-images/JD7348A/5.jpg,https://oleks-netizen.github.io/product-images/JD7348A/4.jpg,https://oleks-netizen.github.io/product-images/JD7348A/6.jpg,https://oleks-netizen.github.io/product-images/JD7348A/7.jpg,https://oleks-netizen.github.io/product-images/JD7348A/9.jpg"
$ws.Range("C17").Value = 7

$ws.Range("A18").Value = "K1266-1bl-black"
$ws.Range("B18").Value = "https://oleks-netizen.github.io/product-images/K1266-1bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K1266-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K1266-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K1266-1bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K1266-1bl-black/6.jpg,https://oleks-netizen.github.io/product-images/K1266-1bl-black/3.jpg"
$ws.Range("C18").Value = 6

$ws.Range("A19").Value = "K1266-2bl-black"
$ws.Range("B19").Value = "https://oleks-netizen.github.io/product-images/K1266-2bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K1266-2bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K1266-2bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K1266-2bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K1266-2bl-black/6.jpg,https://oleks-netizen.github.io/product-images/K1266-2bl-black/4.jpg"
$ws.Range("C19").Value = 6

$ws.Range("A20").Value = "K1426f-black"
$ws.Range("B20").Value = "https://oleks-netizen.github.io/product-images/K1426f-black/1.jpg,https://oleks-netizen.github.io/product-images/K1426f-black/2.jpg,https://oleks-netizen.github.io/product-images/K1426f-black/5.jpg,https://oleks-netizen.github.io/product-images/K1426f-black/4.jpg"
$ws.Range("C20").Value = 4

$ws.Range("A21").Value = "K1428f-black"
$ws.Range("B21").Value = "https://oleks-netizen.github.io/product-images/K1428f-black/2.jpg,https://oleks-netizen.github.io/product-images/K1428f-black/3.jpg,https://oleks-netizen.github.io/product-images/K1428f-black/5.jpg,https://oleks-netizen.github.io/product-images/K1428f-black/4.jpg"
$ws.Range("C21").Value = 4

$ws.Range("A22").Value = "K166030bl-black"
$ws.Range("B22").Value = "https://oleks-netizen.github.io/product-images/K166030bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K166030bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K166030bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K166030bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K166030bl-black/5.jpg,https://oleks-netizen.github.io/product-images/K166030bl-black/6.jpg"
$ws.Range("C22").Value = 6

$ws.Range("A23").Value = "K166317bl-black"
$ws.Range("B23").Value = "https://oleks-netizen.github.io/product-images/K166317bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K166317bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K166317bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K166317bl-black/5.jpg"
$ws.Range("C23").Value = 4

$ws.Range("A24").Value = "K166365bl-black"
$ws.Range("B24").Value = "https://oleks-netizen.github.io/product-images/K166365bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/6.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K166365bl-black/4.jpg"
$ws.Range("C24").Value = 7

$ws.Range("A25").Value = "K16685-1bl-black"
$ws.Range("B25").Value = "https://oleks-netizen.github.io/product-images/K16685-1bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K16685-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K16685-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K16685-1bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K16685-1bl-black/6.jpg,https://oleks-netizen.github.io/product-images/K16685-1bl-black/4.jpg"
$ws.Range("C25").Value = 6

$ws.Range("A26").Value = "K16685-3bl-black"
$ws.Range("B26").Value = "https://oleks-netizen.github.io/product-images/K16685-3bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K16685-3bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K16685-3bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K16685-3bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K16685-3bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K16685-3bl-black/5.jpg"
$ws.Range("C26").Value = 6

$ws.Range("A27").Value = "K19803-1bl-black"
$ws.Range("B27").Value = "https://oleks-netizen.github.io/product-images/K19803-1bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K19803-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K19803-1bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K19803-1bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K19803-1bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K19803-1bl-black/6.jpg"
$ws.Range("C27").Value = 6

$ws.Range("A28").Value = "K19803-1br-brown"
$ws.Range("B28").Value = "https://oleks-netizen.github.io/product-images/K19803-1br-brown/1.jpg,https://oleks-netizen.github.io/product-images/K19803-1br-brown/2.jpg,https://oleks-netizen.github.io/product-images/K19803-1br-brown/2.jpg,https://oleks-netizen.github.io/product-images/K19803-1br-brown/3.jpg,https://oleks-netizen.github.io/product-images/K19803-1br-brown/6.jpg,https://oleks-netizen.github.io/product-images/K19803-1br-brown/4.jpg"
$ws.Range("C28").Value = 6

$ws.Range("A29").Value = "K19803-2bl-black"
$ws.Range("B29").Value = "https://oleks-netizen.github.io/product-images/K19803-2bl-black/1.jpg,https://oleks-netizen.github.io/product-images/K19803-2bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K19803-2bl-black/2.jpg,https://oleks-netizen.github.io/product-images/K19803-2bl-black/3.jpg,https://oleks-netizen.github.io/product-images/K19803-2bl-black/4.jpg,https://oleks-netizen.github.io/product-images/K19803-2bl-black/6.jpg"
$ws.Range("C29").Value = 6

$ws.Range("A30").Value = "LB105 GRY"
$ws.Range("B30").Value = "https://oleks-netizen.github.io/product-images/LB105 GRY/1.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/2.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/3.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/4.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/5.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/6.jpg,https://oleks-netizen.github.io/product-images/LB105 GRY/7.jpg"
$ws.Range("C30").Value = 7

$ws.Range("A31").Value = "LBJX010013"
$ws.Range("B31").Value = "https://oleks-netizen.github.io/product-images/LBJX010013/2.jpg,https://oleks-netizen.github.io/product-images/LBJX010013/3.jpg,https://oleks-netizen.github.io/product-images/LBJX010013/5.jpg,https://oleks-netizen.github.io/product-images/LBJX010013/6.jpg,https://oleks-netizen.github.io/product-images/LBJX010013/4.jpg"
$ws.Range("C31").Value = 5

$ws.Range("A32").Value = "RE-3079-3md"
$ws.Range("B32").Value = "https://oleks-netizen.github.io/product-images/RE-3079-3md/2.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/4.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/11.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/10.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/3.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/5.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/6.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/7.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/8.jpg,https://oleks-netizen.github.io/product-images/RE-3079-3md/9.jpg"
$ws.Range("C32").Value = 10

$ws.Range("A33").Value = "RH-1811-4lx"
$ws.Range("B33").Value = "https://oleks-netizen.github.io/product-images/RH-1811-4lx/1.jpg"
$ws.Range("C33").Value = 1

$ws.Rows(36).Delete()
$ws.Rows(35).Delete()
$ws.Rows(34).Delete()
